# fix: MRR bridge inclusive end-date logic (Diff 256.0 resolved)
# Recomputed Monthly_Data bridge (MRR_total, MRR_start, new_MRR, churned_MRR,
# active_paid_users, ARPU) for rows 3-13 after correcting the end-date
# inclusivity bug in the monthly cohort aggregation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monthly_Data")

# row 3 (month 2024-02)
$ws.Cells.Item(3, 2).Value  = 47841                 # B3  MRR_total
$ws.Cells.Item(3, 7).Value  = 937                   # G3  churned_MRR
$ws.Cells.Item(3, 8).Value  = 242                   # H3  active_paid_users
$ws.Cells.Item(3, 9).Value  = 197.6900826446281     # I3  ARPU

# row 4 (month 2024-03)
$ws.Cells.Item(4, 2).Value  = 79097                 # B4
$ws.Cells.Item(4, 3).Value  = 47841                 # C4  MRR_start
$ws.Cells.Item(4, 4).Value  = 32024                 # D4  new_MRR
$ws.Cells.Item(4, 8).Value  = 391                   # H4
$ws.Cells.Item(4, 9).Value  = 202.2941176470588     # I4

# row 5 (month 2024-04)
$ws.Cells.Item(5, 2).Value  = 107093                # B5
$ws.Cells.Item(5, 3).Value  = 79097                 # C5
$ws.Cells.Item(5, 4).Value  = 28991                 # D5
$ws.Cells.Item(5, 7).Value  = 995                   # G5
$ws.Cells.Item(5, 8).Value  = 534                   # H5
$ws.Cells.Item(5, 9).Value  = 200.5486891385768     # I5

# row 6 (month 2024-05)
$ws.Cells.Item(6, 2).Value  = 130516                # B6
$ws.Cells.Item(6, 3).Value  = 107093                # C6
$ws.Cells.Item(6, 4).Value  = 24488                 # D6
$ws.Cells.Item(6, 7).Value  = 1065                  # G6
$ws.Cells.Item(6, 8).Value  = 635                   # H6
$ws.Cells.Item(6, 9).Value  = 205.5370078740158     # I6

# row 7 (month 2024-06)
$ws.Cells.Item(7, 2).Value  = 151855                # B7
$ws.Cells.Item(7, 3).Value  = 130516                # C7
$ws.Cells.Item(7, 4).Value  = 21880                 # D7
$ws.Cells.Item(7, 7).Value  = 541                   # G7
$ws.Cells.Item(7, 8).Value  = 746                   # H7
$ws.Cells.Item(7, 9).Value  = 203.558981233244      # I7

# row 8 (month 2024-07)
$ws.Cells.Item(8, 2).Value  = 170282                # B8
$ws.Cells.Item(8, 3).Value  = 151855                # C8
$ws.Cells.Item(8, 4).Value  = 19113                 # D8
$ws.Cells.Item(8, 8).Value  = 831                   # H8
$ws.Cells.Item(8, 9).Value  = 204.9121540312876     # I8

# row 9 (month 2024-08)
$ws.Cells.Item(9, 2).Value  = 192220                # B9
$ws.Cells.Item(9, 3).Value  = 170282                # C9
$ws.Cells.Item(9, 4).Value  = 22578                 # D9
$ws.Cells.Item(9, 7).Value  = 640                   # G9
$ws.Cells.Item(9, 8).Value  = 951                   # H9
$ws.Cells.Item(9, 9).Value  = 202.124079915878      # I9

# row 10 (month 2024-09)
$ws.Cells.Item(10, 2).Value = 216536                # B10
$ws.Cells.Item(10, 3).Value = 192220                # C10
$ws.Cells.Item(10, 4).Value = 24886                 # D10
$ws.Cells.Item(10, 7).Value = 570                   # G10
$ws.Cells.Item(10, 8).Value = 1083                  # H10
$ws.Cells.Item(10, 9).Value = 199.9409048938135     # I10

# row 11 (month 2024-10)
$ws.Cells.Item(11, 2).Value = 248355                # B11
$ws.Cells.Item(11, 3).Value = 216536                # C11
$ws.Cells.Item(11, 4).Value = 32645                 # D11
$ws.Cells.Item(11, 8).Value = 1222                  # H11
$ws.Cells.Item(11, 9).Value = 203.2364975450081     # I11

# row 12 (month 2024-11)
$ws.Cells.Item(12, 2).Value = 282968                # B12
$ws.Cells.Item(12, 3).Value = 248355                # C12
$ws.Cells.Item(12, 4).Value = 36190                 # D12
$ws.Cells.Item(12, 7).Value = 1577                  # G12
$ws.Cells.Item(12, 8).Value = 1368                  # H12
$ws.Cells.Item(12, 9).Value = 206.8479532163742     # I12

# row 13 (month 2024-12)
$ws.Cells.Item(13, 2).Value = 312065                # B13
$ws.Cells.Item(13, 3).Value = 282968                # C13
$ws.Cells.Item(13, 4).Value = 30312                 # D13
$ws.Cells.Item(13, 7).Value = 1215                  # G13
$ws.Cells.Item(13, 8).Value = 1488                  # H13
$ws.Cells.Item(13, 9).Value = 209.7211021505376     # I13
